# Refresh the market-data snapshot held in the workbook.
#
# "Neg_Change" gains a new row (IIFL) -> dimension grows from A1:I11 to A1:I12.
# "Pos_Change" loses two rows (PATANJALI/GLENMARK/ABFRL removed, replaced by a
# shorter list ending at IEX) -> dimension shrinks from A1:I19 to A1:I17.
# All OHLC/volume/change% figures for the remaining symbols are refreshed too.

$negData = @(
    @("HINDZINC",516.5,516.75,506.05,506.25,3256458,6882125,-0.5268237644622845,"HINDZINC"),
    @("MARICO",710,726.95,708.4,726.75,1205842,2960078,-0.5926316806516585,"MARICO"),
    @("PRESTIGE",1677.2,1714.9,1677.2,1690,1021953,2176330,-0.5304236949359702,"PRESTIGE"),
    @("360ONE",1174,1179.8,1158,1169.7,915145,2005553,-0.5436944324084181,"360ONE"),
    @("MUTHOOTFIN",3265,3279.8,3252.1,3267.6,234148,529423,-0.5577298304002659,"MUTHOOTFIN"),
    @("NATIONALUM",227,227.75,223.92,226.7,6327250,12613195,-0.4983626273913945,"NATIONALUM"),
    @("AUBANK",801.25,814,795,795,1864887,3804314,-0.5097967728215915,"AUBANK"),
    @("NMDC",76.81,76.85,75.66,75.81,7993203,16267362,-0.5086355734875759,"NMDC"),
    @("ABFRL",83.45,84.19,82.81,83.76,2193997,4595235,-0.5225495540489224,"ABFRL"),
    @("NUVAMA",7220,7300,7200,7230,51527,102908,-0.4992906285225638,"NUVAMA"),
    @("IIFL",508.2,509.8,498,499.55,1142761,2376538,-0.5191488627575069,"IIFL"),
)

$posData = @(
    @("TITAN",3549,3652,3540,3636,1528118,1005045,0.5204473431537892,"TITAN"),
    @("ITC",401,406,399.55,405.45,22962049,16043867,0.4312041479775418,"ITC"),
    @("LODHA",1177,1197,1175,1194,2296704,1575396,0.4578582146964953,"LODHA"),
    @("ADANIGREEN",1055,1067,1052.5,1059,1742078,1161123,0.500338896051495,"ADANIGREEN"),
    @("ZYDUSLIFE",977.3,989,974.75,983.95,832737,563166,0.4786705873579016,"ZYDUSLIFE"),
    @("JSWENERGY",545.75,551.5,545,548.2,2481559,1730184,0.4342746205027905,"JSWENERGY"),
    @("TORNTPHARM",3535,3549.1,3520.9,3531,207979,131171,0.5855562586242385,"TORNTPHARM"),
    @("SOLARINDS",14052,14134,13960,14040,81585,56766,0.4372159391184864,"SOLARINDS"),
    @("PNB",116.65,117.84,115.97,116.13,15188935,10598247,0.433155407682044,"PNB"),
    @("CHOLAFIN",1698,1700.1,1653.2,1664.7,3202080,2004535,0.5974178550137563,"CHOLAFIN"),
    @("JSL",785,812,779.65,808.2,696824,471408,0.478176017377728,"JSL"),
    @("JUBLFOOD",584.7,597,584.55,593.4,1285374,879705,0.4611420874043003,"JUBLFOOD"),
    @("GMRAIRPORT",89.96,91.15,89.79,90.63,6688742,4646050,0.4396620785398349,"GMRAIRPORT"),
    @("LTF",265,273.57,264.92,269.2,11133785,7835400,0.4209593639124997,"LTF"),
    @("NBCC",112,114.63,111.68,112.25,9669421,6608121,0.4632633088891683,"NBCC"),
    @("IEX",135,136.39,133.5,134.97,6846625,4802387,0.4256712339093038,"IEX"),
)

$wb = $excel.ActiveWorkbook

# --- Update "Neg_Change" sheet (formerly rows 2-11, now rows 2-12) ---
$wsNeg = $wb.Worksheets.Item("Neg_Change")
$r = 2
foreach ($row in $negData) {
    for ($c = 0; $c -lt $row.Count; $c++) {
        $wsNeg.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}

# --- Update "Pos_Change" sheet (formerly rows 2-19, now rows 2-17) ---
$wsPos = $wb.Worksheets.Item("Pos_Change")
$r = 2
foreach ($row in $posData) {
    for ($c = 0; $c -lt $row.Count; $c++) {
        $wsPos.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}

# Remove the now-unused trailing rows (old data had rows up to 19; new data ends at row 17)
$wsPos.Range("A18:I19").ClearContents()
